$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    # Force a literal text value (avoids Excel's automatic number/date
    # inference for strings like "2026-02-24"), then drop back to the
    # default "Normal" style so no stray number-format style sticks to
    # the cell.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---- Row 15 -------------------------------------------------------------
$ws.Cells.Item(15, 1).Value = 131289440
$ws.Cells.Item(15, 2).Value = 57881
$ws.Cells.Item(15, 4).Value = "NT"
$ws.Cells.Item(15, 5).Value = 100049
$ws.Cells.Item(15, 6).Value = "Spillkråka"
$ws.Cells.Item(15, 7).Value = "Dryocopus martius"
$ws.Cells.Item(15, 8).Value = "(Linnaeus, 1758)"
Set-TextCell 15 9 ""
Set-TextCell 15 11 ""
Set-TextCell 15 12 ""
$ws.Cells.Item(15, 13).Value = "färska spår"
Set-TextCell 15 14 ""
$ws.Cells.Item(15, 16).Value = "Luvebo 2:1, Ög"
$ws.Cells.Item(15, 17).Value = 567626
$ws.Cells.Item(15, 18).Value = 6509625
$ws.Cells.Item(15, 19).Value = 10
$ws.Cells.Item(15, 20).Value = "Östergötland"
$ws.Cells.Item(15, 21).Value = "Norrköping"
$ws.Cells.Item(15, 22).Value = "Östergötland"
$ws.Cells.Item(15, 23).Value = "Simonstorp"
Set-TextCell 15 25 "2026-02-24"
Set-TextCell 15 27 "2026-02-24"
$ws.Cells.Item(15, 30).Value = $false
$ws.Cells.Item(15, 31).Value = $false
$ws.Cells.Item(15, 33).Value = $false
Set-TextCell 15 46 ""
$ws.Cells.Item(15, 49).Value = "Anette Källman"
$ws.Cells.Item(15, 50).Value = "Anette Källman"
Set-TextCell 15 51 ""

# ---- Row 16 -------------------------------------------------------------
$ws.Cells.Item(16, 1).Value = 131289461
$ws.Cells.Item(16, 2).Value = 57881
$ws.Cells.Item(16, 4).Value = "NT"
$ws.Cells.Item(16, 5).Value = 100049
$ws.Cells.Item(16, 6).Value = "Spillkråka"
$ws.Cells.Item(16, 7).Value = "Dryocopus martius"
$ws.Cells.Item(16, 8).Value = "(Linnaeus, 1758)"
Set-TextCell 16 9 ""
Set-TextCell 16 11 ""
Set-TextCell 16 12 ""
$ws.Cells.Item(16, 13).Value = "gammalt bo"
Set-TextCell 16 14 ""
$ws.Cells.Item(16, 16).Value = "Luvebo 2:1, Ög"
$ws.Cells.Item(16, 17).Value = 567680
$ws.Cells.Item(16, 18).Value = 6509667
$ws.Cells.Item(16, 19).Value = 10
$ws.Cells.Item(16, 20).Value = "Östergötland"
$ws.Cells.Item(16, 21).Value = "Norrköping"
$ws.Cells.Item(16, 22).Value = "Östergötland"
$ws.Cells.Item(16, 23).Value = "Simonstorp"
Set-TextCell 16 25 "2026-02-24"
Set-TextCell 16 27 "2026-02-24"
$ws.Cells.Item(16, 30).Value = $false
$ws.Cells.Item(16, 31).Value = $false
$ws.Cells.Item(16, 33).Value = $false
Set-TextCell 16 46 ""
$ws.Cells.Item(16, 49).Value = "Anette Källman"
$ws.Cells.Item(16, 50).Value = "Anette Källman"
Set-TextCell 16 51 ""
